$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(230, 44304, 1, 3, 122.3990208078335),
    @(231, 44305, 0, 2, 81.59934720522236),
    @(232, 44306, 3, 4, 163.1986944104447),
    @(233, 44307, 1, 5, 203.9983680130559)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy the formatting of the last existing data row (229) down to the new row
    # so column A keeps its centered/bordered date style.
    $ws.Range("A229:D229").Copy()
    $ws.Range("A$r`:D$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = $false
